# Insert a new record row at row 646 in the daily Mango price log
# (Fruta, Vega Central Mapocho de Santiago - Mango), shifting the
# existing rows 646:667 down to 647:668, then fill in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 646:667 down by inserting a new blank row at 646.
$ws.Rows.Item(646).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A646").Value = 9
$ws.Range("B646").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C646").Value = "Metropolitana"
$ws.Range("D646").Value = 45075
$ws.Range("E646").Value = 13
$ws.Range("F646").Value = "Fruta"
$ws.Range("G646").Value = 100108
$ws.Range("H646").Value = "Tropicales y subtropicales"
$ws.Range("I646").Value = 100108002
$ws.Range("J646").Value = "Mango"
$ws.Range("K646").Value = "Sin especificar"
$ws.Range("L646").Value = "Primera"
$ws.Range("M646").Value = 450
$ws.Range("N646").Value = 7500
$ws.Range("O646").Value = 8500
$ws.Range("P646").Value = 7944
$ws.Range("Q646").Value = "`$/bandeja 4 kilos"
$ws.Range("R646").Value = "Perú"
$ws.Range("S646").Value = 1986
$ws.Range("T646").Value = 4
